$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest data row (previously row 2). Excel shifts every row below it
# up by one, which is exactly how A (Datetime) and B (Real_Close) arrive at their
# new values, and how the old Real_Close column lines up under column C.
$ws.Rows.Item(2).Delete()

# Trad_Prediction (C) is refreshed to equal the prior period Real_Close, and
# AI_Prediction (D) is replaced with a newly run simulation series.
$ws.Cells.Item(2, 3).Value = 163.0686950683594
$ws.Cells.Item(2, 4).Value = 155.3970684269742
$ws.Cells.Item(3, 3).Value = 163.3500061035156
$ws.Cells.Item(3, 4).Value = 165.3898519224469
$ws.Cells.Item(4, 3).Value = 163.2700042724609
$ws.Cells.Item(4, 4).Value = 160.7850187702587
$ws.Cells.Item(5, 3).Value = 162.5650024414062
$ws.Cells.Item(5, 4).Value = 167.8253973935246
$ws.Cells.Item(6, 3).Value = 162.8150024414062
$ws.Cells.Item(6, 4).Value = 164.7868737954554
$ws.Cells.Item(7, 3).Value = 162.8500061035156
$ws.Cells.Item(7, 4).Value = 162.6922219202335
$ws.Cells.Item(8, 3).Value = 162.8300018310547
$ws.Cells.Item(8, 4).Value = 162.5305082707869
$ws.Cells.Item(9, 3).Value = 162.4100036621094
$ws.Cells.Item(9, 4).Value = 163.1237049094994
$ws.Cells.Item(10, 3).Value = 162.5966949462891
$ws.Cells.Item(10, 4).Value = 156.0867288086234
$ws.Cells.Item(11, 3).Value = 163.2747039794922
$ws.Cells.Item(11, 4).Value = 156.8594619335315
$ws.Cells.Item(12, 3).Value = 163.4476928710938
$ws.Cells.Item(12, 4).Value = 162.8814893027454
$ws.Cells.Item(13, 3).Value = 163.7299957275391
$ws.Cells.Item(13, 4).Value = 167.8254144098223
$ws.Cells.Item(14, 3).Value = 164.0599975585938
$ws.Cells.Item(14, 4).Value = 158.9572847164069
$ws.Cells.Item(15, 3).Value = 164.1049957275391
$ws.Cells.Item(15, 4).Value = 161.2805195773754
$ws.Cells.Item(16, 3).Value = 167.3598937988281
$ws.Cells.Item(16, 4).Value = 161.4305983826725
$ws.Cells.Item(17, 3).Value = 166.4949951171875
$ws.Cells.Item(17, 4).Value = 168.0552166658724
$ws.Cells.Item(18, 3).Value = 166.4786987304688
$ws.Cells.Item(18, 4).Value = 170.2978857305087
$ws.Cells.Item(19, 3).Value = 166.1822967529297
$ws.Cells.Item(19, 4).Value = 161.0109778602229
$ws.Cells.Item(20, 3).Value = 165.6815032958984
$ws.Cells.Item(20, 4).Value = 171.4147473769869
$ws.Cells.Item(21, 3).Value = 165.0950012207031
$ws.Cells.Item(21, 4).Value = 164.5697456719568
$ws.Cells.Item(22, 3).Value = 164.8800048828125
$ws.Cells.Item(22, 4).Value = 164.2441569051837
$ws.Cells.Item(23, 3).Value = 163.7200012207031
$ws.Cells.Item(23, 4).Value = 158.9066503716691
$ws.Cells.Item(24, 3).Value = 164.7546997070312
$ws.Cells.Item(24, 4).Value = 167.0863752835365
$ws.Cells.Item(25, 3).Value = 164.9700012207031
$ws.Cells.Item(25, 4).Value = 168.5064716028009
$ws.Cells.Item(26, 3).Value = 165.0800933837891
$ws.Cells.Item(26, 4).Value = 169.8187740713424
$ws.Cells.Item(27, 3).Value = 164.8200073242188
$ws.Cells.Item(27, 4).Value = 161.6771507449955
$ws.Cells.Item(28, 3).Value = 164.4033050537109
$ws.Cells.Item(28, 4).Value = 164.5195951904789
$ws.Cells.Item(29, 3).Value = 164.0700073242188
$ws.Cells.Item(29, 4).Value = 163.1512261523607
$ws.Cells.Item(30, 3).Value = 171.5783996582031
$ws.Cells.Item(30, 4).Value = 175.41558571811
$ws.Cells.Item(31, 3).Value = 171.3849945068359
$ws.Cells.Item(31, 4).Value = 173.2348727227712
$ws.Cells.Item(32, 3).Value = 171.3849945068359
$ws.Cells.Item(32, 4).Value = 173.5243658436468
$ws.Cells.Item(33, 3).Value = 170.5249938964844
$ws.Cells.Item(33, 4).Value = 170.356498387303
$ws.Cells.Item(34, 3).Value = 170.9118041992188
$ws.Cells.Item(34, 4).Value = 178.2545713311746
$ws.Cells.Item(35, 3).Value = 170.1999053955078
$ws.Cells.Item(35, 4).Value = 164.3797181212215
